# Apply the "base" showcase macro edit:
#   - add new "base" function outputToCloud(resource)
#   - add new "target" category "text" (with its sole function spellCheck(var,profile,text))
#
# The '#system' sheet stores, per column, the sorted list of values backing
# each named-range / data-validation list used elsewhere in the workbook.
# Column A = list of category ("target") names.
# Column E = the "base" category's function list.
# Column Y (new) = the "text" category's function list; the old Y..AD block
#   (web, webalert, webcookie, ws, ws.async, xml) all shift one column right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Column A ("target" list): insert "text" at A25, push A25:A30 -> A26:A31
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------------
# 2) Column E ("base" function list): insert "outputToCloud(resource)" at
#    E21, push E21:E37 -> E22:E38
# ---------------------------------------------------------------------
for ($r = 37; $r -ge 21; $r--) {
    $v = $ws.Cells.Item($r, 5).Value2()
    $ws.Cells.Item($r + 1, 5).Value = $v
}
$ws.Cells.Item(21, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 3) Insert a brand-new column at Y (25): shifts Y..AD one column right
#    (web, webalert, webcookie, ws, ws.async, xml all move right one column).
# ---------------------------------------------------------------------
$ws.Columns.Item(25).Insert()

# New column Y becomes the "text" category list: header + sole function.
$ws.Cells.Item(1, 25).Value = "text"
$ws.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 4) Fix up the workbook-level defined names so they describe the new
#    layout (the engine does not auto-adjust these on insert/shift).
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")

Write-Output "done"
